# Generate Report for Handoff
# Adds a new localization-status row (for e2e\ed771909-4059-46b2-9328-21c17446c40c.md)
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$fileName      = "ed771909-4059-46b2-9328-21c17446c40c.md"
$pathAndName   = "e2e\ed771909-4059-46b2-9328-21c17446c40c.md"
$status        = "Ready for handoff"
$zhHandoffFile = "ed771909-4059-46b2-9328-21c17446c40c.734f1fecb52a5e7e20239177d9868deb17caf766.zh-cn.xlf"
$deHandoffFile = "ed771909-4059-46b2-9328-21c17446c40c.734f1fecb52a5e7e20239177d9868deb17caf766.de-de.xlf"
$zhHandoffDate = "2016-08-24 10:51:28"
$deHandoffDate = "2016-08-24 10:51:32"
$noHandback    = "0001-01-01 00:00:00"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f49c498468ba7a60a83ba4fa2cc28cc1129989ba/e2e/$fileName"
$repoZhCn = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1cbc96e46767ba0272cad2acfb73043449731ebe/e2e/$fileName"
$repoDeDe = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7e044279feb815ff66d9a9b25c676c308eb51c5e/e2e/$fileName"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $fileName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $deHandoffDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $repoBase, "", "", $pathAndName) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = $zhHandoffFile
$wsZhCn.Range("H3").Value = $zhHandoffDate
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = $noHandback
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $repoZhCn, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = $deHandoffFile
$wsDeDe.Range("H3").Value = $deHandoffDate
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = $noHandback
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $repoDeDe, "", "", $fileName) | Out-Null
